# "Generate Report for Archive"
# Refresh the localization-status report: the two "Ready for handoff" test
# rows have moved on to "In Translation", and the Status columns are
# re-sized to fit the (now shorter) status text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E, F) -------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Narrower column now that the status text is shorter than
# "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: Status column (C) ----------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
